$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40: log a new completed task - "Able to swipe to hide the article."
$ws.Range("B40").Value = 2
$ws.Range("C40").Value = "Able to swipe to hide the article."

# Extend the date list with three more days (rows 42-44), continuing from A41
$base = $ws.Range("A41").Value2
$ws.Range("A42").Value = $base + 1
$ws.Range("A43").Value = $base + 2
$ws.Range("A44").Value = $base + 3

# Match the date formatting/style used by the preceding date cells
$ws.Range("A41").Copy() | Out-Null
$ws.Range("A42:A44").PasteSpecial(-4122) | Out-Null

# Update the view - scroll back to the top and select the newly added rows
$ws.Activate() | Out-Null
$ws.Range("A38:A44").Select() | Out-Null
